# Weekly update: prepend newly reported price records for "Perejil" /
# Vega Modelo de Temuco, pushing the existing history down.
#
# Net effect (matches the target diff):
#   - Two brand-new rows are inserted at the top of the data block (rows 290-291),
#     moving the former rows 290-325 down to 292-327.
#   - One more brand-new row is inserted further down (ends up at row 315),
#     moving the former row 313 (now at 315) and everything after it down one more,
#     to 316-328.
#   - dimension grows from A1:R325 to A1:R328.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow {
    param($row, $date, $volumen, $pmin, $pmax, $pprom, $origen, $precioKg)

    $ws.Cells.Item($row, 1).Value = 10
    $ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($row, 3).Value = "La Araucanía"
    $ws.Cells.Item($row, 4).Value = $date
    $ws.Cells.Item($row, 5).Value = 9
    $ws.Cells.Item($row, 6).Value = 100112044
    $ws.Cells.Item($row, 7).Value = "Perejil"
    $ws.Cells.Item($row, 8).Value = "Sin especificar"
    $ws.Cells.Item($row, 9).Value = "Primera"
    $ws.Cells.Item($row, 10).Value = $volumen
    $ws.Cells.Item($row, 11).Value = $pmin
    $ws.Cells.Item($row, 12).Value = $pmax
    $ws.Cells.Item($row, 13).Value = $pprom
    $ws.Cells.Item($row, 14).Value = "$/docena de atados (3 kilos)"
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $precioKg
    $ws.Cells.Item($row, 17).Value = 3
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}

# --- Insert the two newest records at the top of the block (rows 290-291) ---
$ws.Rows.Item(290).Insert()
$ws.Rows.Item(290).Insert()

Set-DataRow 290 44748 35 4000 4000 4000 "Provincia de Cautín" 1333
Set-DataRow 291 44748 35 3666 3666 3666 "Región Metropolitana" 1222

# --- Insert one more record further down the history (lands on row 315) ---
$ws.Rows.Item(315).Insert()

Set-DataRow 315 44747 60 4000 5000 4417 "Provincia de Cautín" 1472
